# Daily update at 8 AM UTC
# Appends the next day's row of data to the "Wins Over Time" tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find the next empty row right after the existing data (row 36 -> new row 37)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# New day's values
$ws.Cells.Item($newRow, 1).Value = 45986
$ws.Cells.Item($newRow, 2).Value = 83
$ws.Cells.Item($newRow, 3).Value = 90
$ws.Cells.Item($newRow, 4).Value = 90

# Match the date-time number formatting used by the rest of column A
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($lastRow, 1).NumberFormat
